# Refresh the cryptos price/volume table with the latest scraped values.
# Note: the Price column stores plain-text numbers (e.g. "312.23"); a leading
# apostrophe ('') forces Excel to keep them as text instead of coercing them
# into floating point numbers (which would mangle values like "1.110" or
# "0.06501" by dropping significant trailing/leading zeros).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.999.98'
$ws.Range("E2").Value = '  -0.33%  '
$ws.Range("D3").Value = '1.858.32'
$ws.Range("E3").Value = '  -0.94%  '
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").Value = '''312.23'
$ws.Range("E5").Value = '  -0.39%  '
$ws.Range("E6").Value = '  +0.06%  '
$ws.Range("D7").Value = '''0.5136'
$ws.Range("E7").Value = '  +1.42%  '
$ws.Range("D8").Value = '''0.3833'
$ws.Range("E8").Value = '  -0.30%  '
$ws.Range("D9").Value = '''0.08243'
$ws.Range("E9").Value = '  -5.22%  '
$ws.Range("D10").Value = '''1.111'
$ws.Range("E10").Value = '  -0.75%  '
$ws.Range("D11").Value = '''41.47'
$ws.Range("E11").Value = '  -0.14%  '
$ws.Range("D12").Value = '''6.195'
$ws.Range("E12").Value = '  -2.09%  '
$ws.Range("E13").Value = '  -0.64%  '
$ws.Range("D14").Value = '1.861.32'
$ws.Range("E14").Value = '  -1.07%  '
$ws.Range("D15").Value = '''7.260'
$ws.Range("E15").Value = '  +1.13%  '
$ws.Range("E16").Value = '  -0.04%  '
$ws.Range("E17").Value = '  -0.66%  '
$ws.Range("D18").Value = '''90.50'
$ws.Range("E18").Value = '  -0.58%  '
$ws.Range("E19").Value = '  +0.56%  '
$ws.Range("D20").Value = '''17.67'
$ws.Range("E20").Value = '  -2.56%  '
$ws.Range("E21").Value = '  +0.15%  '
$ws.Range("E22").Value = '  -1.32%  '
$ws.Range("D23").Value = '28.029.46'
$ws.Range("E23").Value = '  -0.45%  '
$ws.Range("E24").Value = '  -3.02%  '
$ws.Range("D25").Value = '''2.256'
$ws.Range("E25").Value = '  -0.84%  '
$ws.Range("D26").Value = '2.072.27'
$ws.Range("E26").Value = '  -1.15%  '
$ws.Range("D27").Value = '''2.506'
$ws.Range("E27").Value = '  -1.74%  '
$ws.Range("D28").Value = '''157.29'
$ws.Range("E28").Value = '  +0.12%  '
$ws.Range("D29").Value = '''20.46'
$ws.Range("E29").Value = '  -1.50%  '
$ws.Range("D30").Value = '''124.64'
$ws.Range("E30").Value = '  -1.64%  '
$ws.Range("D31").Value = '''0.1066'
$ws.Range("E31").Value = '  +1.34%  '
$ws.Range("D32").Value = '''1.031'
$ws.Range("E32").Value = '  -2.60%  '
$ws.Range("D33").Value = '''5.902'
$ws.Range("E33").Value = '  +5.63%  '
$ws.Range("D34").Value = '''3.588'
$ws.Range("D35").Value = '''9.425'
$ws.Range("E35").Value = '  -2.31%  '
$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").Value = '''0.06501'
$ws.Range("E36").Value = '  -0.70%  '
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").Value = '''0.02410'
$ws.Range("E37").Value = '  -1.61%  '
$ws.Range("D38").Value = '''0.2176'
$ws.Range("E38").Value = '  +0.27%  '
$ws.Range("D39").Value = '''0.6541'
$ws.Range("E39").Value = '  +2.66%  '
$ws.Range("D40").Value = '''1.194'
$ws.Range("E40").Value = '  -1.05%  '
$ws.Range("D41").Value = '''4.988'
$ws.Range("E41").Value = '  +1.75%  '
$ws.Range("E42").Value = '  -2.36%  '
$ws.Range("D43").Value = '''11.16'
$ws.Range("E43").Value = '  -3.01%  '
$ws.Range("D44").Value = '''0.6146'
$ws.Range("E44").Value = '  +2.71%  '
$ws.Range("D45").Value = '''12.98'
$ws.Range("E45").Value = '  -2.15%  '
$ws.Range("D46").Value = '''1.279'
$ws.Range("E46").Value = '  +0.13%  '
$ws.Range("D47").Value = '''3.672'
$ws.Range("E47").Value = '  -0.01%  '
$ws.Range("D48").Value = '''2.007'
$ws.Range("E48").Value = '  +0.84%  '
$ws.Range("E49").Value = '  -1.39%  '
$ws.Range("D50").Value = '''120.92'
$ws.Range("E50").Value = '  +0.01%  '
$ws.Range("D51").Value = '''78.12'
$ws.Range("E51").Value = '  -2.04%  '
